$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date.
$ws.Name = "Through 2022-11-27"

# Update the header label for the current-year column.
$ws.Range("I1").Value = "2022 (through 11-27)"

# Update the November and Total figures for the 2022 column.
$ws.Range("I12").Value = 104
$ws.Range("I14").Value = 1502
